{"js": "// Replace the date line and every two-digit multiplication answer cell.\n// Each old value is unique in the document, so a simple search+replace\n// (matchCase, no wildcards) for each pair is safe and unambiguous.\nconst replacements = [\n  [\"2023-12-07 Thursday\", \"2023-12-08 Friday\"],\n  [\"71\u00d732=2272\", \"99\u00d775=7425\"],\n  [\"97\u00d748=4656\", \"97\u00d778=7566\"],\n  [\"49\u00d750=2450\", \"11\u00d738=418\"],\n  [\"69\u00d723=1587\", \"94\u00d793=8742\"],\n  [\"42\u00d752=2184\", \"24\u00d748=1152\"],\n  [\"63\u00d772=4536\", \"51\u00d790=4590\"],\n  [\"81\u00d715=1215\", \"38\u00d788=3344\"],\n  [\"42\u00d751=2142\", \"27\u00d719=513\"],\n  [\"80\u00d756=4480\", \"31\u00d731=961\"],\n  [\"29\u00d750=1450\", \"41\u00d736=1476\"],\n  [\"37\u00d723=851\", \"70\u00d750=3500\"],\n  [\"38\u00d744=1672\", \"80\u00d779=6320\"],\n  [\"20\u00d769=1380\", \"20\u00d750=1000\"],\n  [\"83\u00d734=2822\", \"28\u00d765=1820\"],\n  [\"89\u00d794=8366\", \"24\u00d713=312\"],\n  [\"74\u00d789=6586\", \"17\u00d740=680\"],\n  [\"27\u00d792=2484\", \"44\u00d740=1760\"],\n  [\"18\u00d718=324\", \"37\u00d778=2886\"],\n  [\"77\u00d752=4004\", \"71\u00d743=3053\"],\n  [\"22\u00d748=1056\", \"71\u00d778=5538\"],\n  [\"99\u00d733=3267\", \"76\u00d726=1976\"],\n  [\"51\u00d796=4896\", \"58\u00d731=1798\"],\n  [\"45\u00d762=2790\", \"99\u00d713=1287\"],\n  [\"58\u00d773=4234\", \"29\u00d764=1856\"],\n  [\"56\u00d773=4088\", \"23\u00d714=322\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every two-digit multiplication answer.\n# Every \"old\" value below is unique in the document, so Find.Execute with\n# Replace:=2 (wdReplaceAll) safely retargets exactly one run each, with no\n# risk of touching an unrelated cell.\n$replacements = @(\n    @(\"2023-12-07 Thursday\", \"2023-12-08 Friday\"),\n    @(\"71\u00d732=2272\", \"99\u00d775=7425\"),\n    @(\"97\u00d748=4656\", \"97\u00d778=7566\"),\n    @(\"49\u00d750=2450\", \"11\u00d738=418\"),\n    @(\"69\u00d723=1587\", \"94\u00d793=8742\"),\n    @(\"42\u00d752=2184\", \"24\u00d748=1152\"),\n    @(\"63\u00d772=4536\", \"51\u00d790=4590\"),\n    @(\"81\u00d715=1215\", \"38\u00d788=3344\"),\n    @(\"42\u00d751=2142\", \"27\u00d719=513\"),\n    @(\"80\u00d756=4480\", \"31\u00d731=961\"),\n    @(\"29\u00d750=1450\", \"41\u00d736=1476\"),\n    @(\"37\u00d723=851\", \"70\u00d750=3500\"),\n    @(\"38\u00d744=1672\", \"80\u00d779=6320\"),\n    @(\"20\u00d769=1380\", \"20\u00d750=1000\"),\n    @(\"83\u00d734=2822\", \"28\u00d765=1820\"),\n    @(\"89\u00d794=8366\", \"24\u00d713=312\"),\n    @(\"74\u00d789=6586\", \"17\u00d740=680\"),\n    @(\"27\u00d792=2484\", \"44\u00d740=1760\"),\n    @(\"18\u00d718=324\", \"37\u00d778=2886\"),\n    @(\"77\u00d752=4004\", \"71\u00d743=3053\"),\n    @(\"22\u00d748=1056\", \"71\u00d778=5538\"),\n    @(\"99\u00d733=3267\", \"76\u00d726=1976\"),\n    @(\"51\u00d796=4896\", \"58\u00d731=1798\"),\n    @(\"45\u00d762=2790\", \"99\u00d713=1287\"),\n    @(\"58\u00d773=4234\", \"29\u00d764=1856\"),\n    @(\"56\u00d773=4088\", \"23\u00d714=322\")\n)\n\n$d = $word.ActiveDocument\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n    if (-not $found) {\n        Write-Output \"NOT FOUND: $oldText\"\n    }\n}\n\nWrite-Output \"done\"\n"}
